$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 9.563969737553307
$ws.Range("F2").Value = 1.733275416061598
$ws.Range("G2").Value = 0.2315502023843321

$ws.Range("E3").Value = 9.572371048910323
$ws.Range("F3").Value = 1.732933490523181
$ws.Range("G3").Value = 0.2308751703433318

$ws.Range("E4").Value = 9.511271359837897
$ws.Range("F4").Value = 1.731058709006017
$ws.Range("G4").Value = 0.2357844334412406

$ws.Range("E5").Value = 9.775052123024182
$ws.Range("F5").Value = 1.735423216370861
$ws.Range("G5").Value = 0.214590067540072

$ws.Range("E6").Value = 9.718653014423852
$ws.Range("F6").Value = 1.740448094675088
$ws.Range("G6").Value = 0.2191216464533189

$ws.Range("E7").Value = 9.617042860501286
$ws.Range("F7").Value = 1.72825230287906
$ws.Range("G7").Value = 0.2272858611424243

$ws.Range("E8").Value = 9.927341565128671
$ws.Range("F8").Value = 1.765927858880759
$ws.Range("G8").Value = 0.2023538524353037

$ws.Range("E9").Value = 9.93508201371888
$ws.Range("F9").Value = 1.762363969788354
$ws.Range("G9").Value = 0.2017319196693262

$ws.Range("E10").Value = 9.781958098876542
$ws.Range("F10").Value = 1.759687280446587
$ws.Range("G10").Value = 0.2140351833349028

$ws.Range("E11").Value = 9.529296916154653
$ws.Range("F11").Value = 1.733781476870487
$ws.Range("G11").Value = 0.2343361085841331

$ws.Range("E12").Value = 9.499849756817117
$ws.Range("F12").Value = 1.731048666847584
$ws.Range("G12").Value = 0.2367021411265023

$ws.Range("E13").Value = 9.525651127705164
$ws.Range("F13").Value = 1.73579260385449
$ws.Range("G13").Value = 0.2346290419029369

$ws.Range("E14").Value = 9.661773719787583
$ws.Range("F14").Value = 1.723711849993406
$ws.Range("G14").Value = 0.2236918075528609

$ws.Range("E15").Value = 9.665319956967608
$ws.Range("F15").Value = 1.732887104214528
$ws.Range("G15").Value = 0.223406873020646

$ws.Range("E16").Value = 9.620376603461862
$ws.Range("F16").Value = 1.738409280039439
$ws.Range("G16").Value = 0.2270180001836744

$ws.Range("E17").Value = 9.861951114293106
$ws.Range("F17").Value = 1.764887508634871
$ws.Range("G17").Value = 0.2076078714346823

$ws.Range("E18").Value = 9.803022144227004
$ws.Range("F18").Value = 1.75953435710057
$ws.Range("G18").Value = 0.2123427207037244

$ws.Range("E19").Value = 9.703578916860224
$ws.Range("F19").Value = 1.762813898051566
$ws.Range("G19").Value = 0.2203328262813484

$ws.Range("E20").Value = 9.50845594215739
$ws.Range("F20").Value = 1.732545153112407
$ws.Range("G20").Value = 0.2360106477859278

$ws.Range("E21").Value = 9.500862350750998
$ws.Range("F21").Value = 1.73262482632591
$ws.Range("G21").Value = 0.2366207808101369

$ws.Range("E22").Value = 9.508466165437833
$ws.Range("F22").Value = 1.733674288314393
$ws.Range("G22").Value = 0.2360098263615609

$ws.Range("E23").Value = 9.663349813922052
$ws.Range("F23").Value = 1.733194590215343
$ws.Range("G23").Value = 0.2235651708892273

$ws.Range("E24").Value = 9.591707081882296
$ws.Range("F24").Value = 1.729261980078666
$ws.Range("G24").Value = 0.2293215507657146

$ws.Range("E25").Value = 9.616078774853047
$ws.Range("F25").Value = 1.744996980839411
$ws.Range("G25").Value = 0.2273633238949846

$ws.Range("E26").Value = 9.857946398417514
$ws.Range("F26").Value = 1.762540254198192
$ws.Range("G26").Value = 0.2079296440028265

$ws.Range("E27").Value = 9.747185518932783
$ws.Range("F27").Value = 1.75975272296495
$ws.Range("G27").Value = 0.2168291049755616

$ws.Range("E28").Value = 9.688587864090525
$ws.Range("F28").Value = 1.763821181764031
$ws.Range("G28").Value = 0.221537333591915
